# Rename the AHB-Diff column headers from the generic "_old"/"_new" suffixes
# to the concrete format-version suffixes "_FV2410"/"_FV2504", then turn the
# used range into a real Excel Table (with AutoFilter) and freeze the header
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- 1. Rename headers (row 1) -------------------------------------------

$oldToFv2410 = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
}

foreach ($addr in $oldToFv2410.Keys) {
    $ws.Range($addr).Value = $oldToFv2410[$addr]
}

$newToFv2504 = @{
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $newToFv2504.Keys) {
    $ws.Range($addr).Value = $newToFv2504[$addr]
}

# --- 2. Convert the used range into an Excel Table with an AutoFilter ----

$tableRange = $ws.Range("A1:U66")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
